$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.092.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.908.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5120"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.14"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3008"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06812"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.909.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.49%  "
$ws.Range("E12").Value = "  +3.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07333"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7007"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.916"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.084.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008232"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +11.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.154.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.840"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.764"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.266"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "148.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "135.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.30%  "
$ws.Range("E28").Value = "  +3.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.004"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.401"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.283"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08825"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.008"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05064"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("E35").Value = "  +4.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7193"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.690"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.807"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9633"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01699"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.169"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4318"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9989"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.622"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1281"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05741"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "33.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.443"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3821"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.46%  "
